# Applies the cryptos list update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '63.369.29'
Set-TextValue $ws.Range("E2") '  -1.06%  '
Set-TextValue $ws.Range("D3") '2.716.04'
Set-TextValue $ws.Range("E3") '  -1.48%  '
Set-TextValue $ws.Range("E4") '  -0.12%  '
Set-TextValue $ws.Range("D5") '558.64'
Set-TextValue $ws.Range("E5") '  -3.13%  '
Set-TextValue $ws.Range("D6") '157.16'
Set-TextValue $ws.Range("E6") '  -1.20%  '
Set-TextValue $ws.Range("E7") '  +0.01%  '
Set-TextValue $ws.Range("E8") '  -1.71%  '
Set-TextValue $ws.Range("E9") '  -3.10%  '
Set-TextValue $ws.Range("E10") '  -0.10%  '
Set-TextValue $ws.Range("D11") '5.59'
Set-TextValue $ws.Range("E11") '  -4.11%  '
Set-TextValue $ws.Range("E12") '  -3.83%  '
Set-TextValue $ws.Range("D13") '3.194.86'
Set-TextValue $ws.Range("E13") '  -1.64%  '
Set-TextValue $ws.Range("D14") '26.46'
Set-TextValue $ws.Range("E14") '  -1.89%  '
Set-TextValue $ws.Range("D15") '63.237.74'
Set-TextValue $ws.Range("E15") '  -0.70%  '
Set-TextValue $ws.Range("D16") '0.0000146'
Set-TextValue $ws.Range("E16") '  -3.79%  '
Set-TextValue $ws.Range("D17") '2.715.84'
Set-TextValue $ws.Range("E17") '  -1.72%  '
Set-TextValue $ws.Range("D18") '12.16'
Set-TextValue $ws.Range("E18") '  -0.44%  '
Set-TextValue $ws.Range("E19") '  -4.40%  '
Set-TextValue $ws.Range("D20") '350.20'
Set-TextValue $ws.Range("E20") '  -2.18%  '
Set-TextValue $ws.Range("E21") '  -4.84%  '
Set-TextValue $ws.Range("E23") '  -4.28%  '
Set-TextValue $ws.Range("D24") '64.33'
Set-TextValue $ws.Range("E24") '  -1.91%  '
Set-TextValue $ws.Range("E25") '  -1.08%  '
Set-TextValue $ws.Range("E26") '  +0.15%  '
Set-TextValue $ws.Range("E27") '  -4.79%  '
Set-TextValue $ws.Range("D28") '0.0₃0880'
Set-TextValue $ws.Range("E28") '  -3.51%  '
Set-TextValue $ws.Range("E29") '  +10.12%  '
Set-TextValue $ws.Range("E30") '  -0.49%  '
Set-TextValue $ws.Range("D31") '7.12'
Set-TextValue $ws.Range("E31") '  -3.07%  '
Set-TextValue $ws.Range("D32") '165.81'
Set-TextValue $ws.Range("E32") '  -1.79%  '
Set-TextValue $ws.Range("E33") '  -0.67%  '
Set-TextValue $ws.Range("D35") '19.81'
Set-TextValue $ws.Range("E35") '  -1.98%  '
Set-TextValue $ws.Range("E36") '  -2.79%  '
Set-TextValue $ws.Range("E37") '  -2.11%  '
Set-TextValue $ws.Range("D38") '344.26'
Set-TextValue $ws.Range("E38") '  -1.79%  '
Set-TextValue $ws.Range("D39") '0.954'
Set-TextValue $ws.Range("E39") '  -4.98%  '
Set-TextValue $ws.Range("D40") '6.01'
Set-TextValue $ws.Range("E40") '  -4.99%  '
Set-TextValue $ws.Range("E41") '  -3.90%  '
Set-TextValue $ws.Range("D42") '38.19'
Set-TextValue $ws.Range("E42") '  -2.57%  '
Set-TextValue $ws.Range("B43") 'EnergySwap'
Set-TextValue $ws.Range("C43") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range("D43") '20.76'
Set-TextValue $ws.Range("E43") '  -3.78%  '
Set-TextValue $ws.Range("B44") 'InjectiveProtocol'
Set-TextValue $ws.Range("C44") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range("D44") '21.22'
Set-TextValue $ws.Range("E44") '  -3.17%  '
Set-TextValue $ws.Range("B45") 'Hedera'
Set-TextValue $ws.Range("C45") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D45") '0.0570'
Set-TextValue $ws.Range("E45") '  -3.59%  '
Set-TextValue $ws.Range("B46") 'Mantle'
Set-TextValue $ws.Range("C46") 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range("D46") '0.625'
Set-TextValue $ws.Range("E46") '  -1.31%  '
Set-TextValue $ws.Range("E47") '  -0.03%  '
Set-TextValue $ws.Range("D48") '131.68'
Set-TextValue $ws.Range("E48") '  -4.07%  '
Set-TextValue $ws.Range("E49") '  -3.34%  '
Set-TextValue $ws.Range("D50") '11.05'
Set-TextValue $ws.Range("E50") '  +0.14%  '
Set-TextValue $ws.Range("D51") '0.0244'
Set-TextValue $ws.Range("E51") '  -4.53%  '
